$d = $word.ActiveDocument

$replacements = @(
    @("37+19=56", "34+52=86"),
    @("98-2=96", "42+47=89"),
    @("31-26=5", "11-4=7"),
    @("93-77=16", "56+3=59"),
    @("64+18=82", "51-44=7"),
    @("11+81=92", "29-6=23"),
    @("11+55=66", "8+33=41"),
    @("44-22=22", "99-78=21"),
    @("0+88=88", "11+61=72"),
    @("11+52=63", "37+39=76"),
    @("35+30=65", "74-32=42"),
    @("23+51=74", "18+21=39"),
    @("81+9=90", "56-8=48"),
    @("35+10=45", "95-84=11"),
    @("18+34=52", "20+29=49"),
    @("65-52=13", "81-32=49"),
    @("86-66=20", "74+13=87"),
    @("63-32=31", "77-11=66"),
    @("67+10=77", "23+66=89"),
    @("67-59=8", "10+51=61"),
    @("60+32=92", "33-11=22"),
    @("63-16=47", "46-41=5"),
    @("64-25=39", "76-41=35"),
    @("94-63=31", "19+16=35"),
    @("78-43=35", "98-25=73"),
    @("83-6=77", "51+48=99"),
    @("64+2=66", "51+47=98"),
    @("11+58=69", "12-10=2"),
    @("22+75=97", "33+64=97"),
    @("59+39=98", "20-0=20"),
    @("88-57=31", "13+17=30"),
    @("30+52=82", "44+23=67"),
    @("3+19=22", "3+28=31"),
    @("52-43=9", "14-12=2"),
    @("55+3=58", "49-47=2"),
    @("13+51=64", "22+69=91"),
    @("60+38=98", "39-26=13"),
    @("48-42=6", "43+4=47"),
    @("22+36=58", "17+60=77"),
    @("20+12=32", "19+38=57"),
    @("95-19=76", "6+48=54"),
    @("92+4=96", "32+60=92"),
    @("33-12=21", "0+44=44"),
    @("82-64=18", "37+32=69"),
    @("78-64=14", "54+30=84"),
    @("42+43=85", "50+35=85"),
    @("88-11=77", "82-52=30"),
    @("18-7=11", "1+78=79"),
    @("2+32=34", "40+50=90"),
    @("0+75=75", "56+24=80"),
    @("44+33=77", "87+11=98"),
    @("80-10=70", "61-10=51"),
    @("90-40=50", "97-77=20"),
    @("65-26=39", "85+5=90"),
    @("90-31=59", "81-80=1"),
    @("36-6=30", "29+40=69"),
    @("74+19=93", "6+51=57"),
    @("61+28=89", "71-6=65"),
    @("64+3=67", "7+76=83"),
    @("80+4=84", "7+11=18"),
    @("66-23=43", "68-31=37"),
    @("3+10=13", "71-16=55"),
    @("97-52=45", "1+97=98"),
    @("84-61=23", "96-27=69"),
    @("54-32=22", "49-47=2"),
    @("26-17=9", "63-52=11"),
    @("66-25=41", "54-29=25"),
    @("77+9=86", "58+7=65"),
    @("43+20=63", "92-2=90"),
    @("39+23=62", "56-41=15"),
    @("26-15=11", "18+19=37"),
    @("84-43=41", "97-46=51"),
    @("42-2=40", "44+42=86"),
    @("7+16=23", "38+36=74"),
    @("6+69=75", "26+7=33"),
    @("59+10=69", "50-43=7"),
    @("0+20=20", "40+50=90"),
    @("36-16=20", "43+32=75"),
    @("89-15=74", "93-13=80"),
    @("8+73=81", "71+5=76"),
    @("99-14=85", "91-21=70"),
    @("69-47=22", "64-26=38"),
    @("8+47=55", "52-41=11"),
    @("41+38=79", "48+31=79"),
    @("99-58=41", "88-85=3"),
    @("88-31=57", "34+40=74"),
    @("38+30=68", "41+17=58"),
    @("60+26=86", "49+4=53"),
    @("65-28=37", "20-11=9"),
    @("31+3=34", "34-9=25"),
    @("47+1=48", "26+72=98"),
    @("50-46=4", "74+5=79"),
    @("28+16=44", "46+10=56"),
    @("73-18=55", "62+0=62"),
    @("40-7=33", "90-62=28"),
    @("12-3=9", "62-48=14"),
    @("76-4=72", "94+3=97"),
    @("19+69=88", "29+45=74"),
    @("43-31=12", "60-44=16"),
    @("88-72=16", "46+39=85"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
